# "Final comit for the night"
# Adds three new slides to the "3.0 FileIO" deck, using the same
# "Title and Content" layout (CustomLayout of slide 2) as the other
# content slides:
#   1. "Warning"           -> inserted right after "Why do we use it?"
#                              (before the StreamReader / StreamWriter pair)
#   2. "Many Other Ways"   -> inserted right after the StreamWriter slide
#   3. "Lab"               -> appended at the very end of the deck
#
# The slide-id allocation order below (Many Other Ways, then Warning,
# then Lab) reproduces the same p:sldId numbering (262, 259, 260, 261,
# 263) that PowerPoint produced for the original edit.

$p = $ppt.ActivePresentation
$layout = $p.Slides.Item(2).CustomLayout

# --- "Many Other Ways" slide (ends up after StreamWriter, slide 7) ---
$sMany = $p.Slides.AddSlide(6, $layout)
$sMany.Shapes.Item(1).Name = "Title 1"
$sMany.Shapes.Item(2).Name = "Text Placeholder 2"
$sMany.Shapes.Item(1).TextFrame.TextRange.Text = "Many Other Ways"

$manyBody = $sMany.Shapes.Item(2).TextFrame.TextRange
$manyUrl = "https://docs.microsoft.com/en-us/dotnet/standard/io/"
$manyBody.Text = $manyUrl + "`r" + `
    "For the scope of this class we will only be talking about a small portion of .NETs FileIO." + "`r" + `
    "Thankfully, we can do a ton of powerful stuff with just the few bits we will discuss."
$manyBody.Characters(1, $manyUrl.Length).ActionSettings.Item(1).Hyperlink.Address = $manyUrl

# --- "Warning" slide (ends up before StreamReader, slide 4) ---
$sWarn = $p.Slides.AddSlide(4, $layout)
$sWarn.Shapes.Item(1).Name = "Title 1"
$sWarn.Shapes.Item(2).Name = "Text Placeholder 2"
$sWarn.Shapes.Item(1).TextFrame.TextRange.Text = "Warning"
$sWarn.Shapes.Item(2).TextFrame.TextRange.Text = "The most common error that you will encounter says that…" + "`r" + `
    "“blah blah does not exist in current context”" + "`r" + `
    "" + "`r" + `
    "Using system.io"

# --- "Lab" slide (appended at the end, slide 8) ---
$sLab = $p.Slides.AddSlide($p.Slides.Count + 1, $layout)
$sLab.Shapes.Item(1).Name = "Title 1"
$sLab.Shapes.Item(2).Name = "Text Placeholder 2"
$sLab.Shapes.Item(1).TextFrame.TextRange.Text = "Lab"
$sLab.Shapes.Item(2).TextFrame.TextRange.Text = "On canvas there is a zip file containing a fairly interesting folder structure with a single file hidden inside." + "`r" + `
    "Your objective is to use file io and recursion to find the path of the file."
